$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update O2 from "Visible" to "Hidden"
$ws.Range("O2").Value = "Hidden"

# Add a new row 7 with project data
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "test"
$ws.Range("C7").Value = "test"
$ws.Range("D7").Value = "2-ROOM"
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = "2-ROOM"
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 44511
$ws.Range("K7").Value = 44876
$ws.Range("L7").Value = "T8765432F"
$ws.Range("M7").Value = 10
$ws.Range("O7").Value = "Visible"
